# Update "想去人数" (column F) values on several sheets, to match the
# regenerated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 203
$ws1.Range("F6").Value  = 1272
$ws1.Range("F9").Value  = 228
$ws1.Range("F15").Value = 215
$ws1.Range("F16").Value = 1567
$ws1.Range("F18").Value = 246
$ws1.Range("F19").Value = 369
$ws1.Range("F21").Value = 880
$ws1.Range("F25").Value = 2714
$ws1.Range("F26").Value = 1497
$ws1.Range("F27").Value = 72
$ws1.Range("F28").Value = 69
$ws1.Range("F29").Value = 499
$ws1.Range("F30").Value = 826
$ws1.Range("F31").Value = 1420
$ws1.Range("F33").Value = 1489
$ws1.Range("F36").Value = 804
$ws1.Range("F37").Value = 690
$ws1.Range("F38").Value = 711
$ws1.Range("F39").Value = 910
$ws1.Range("F41").Value = 271

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 6
$ws2.Range("F12").Value = 2
$ws2.Range("F15").Value = 714
$ws2.Range("F23").Value = 25
$ws2.Range("F25").Value = 40

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 203
$ws4.Range("F10").Value = 1272
$ws4.Range("F13").Value = 228
$ws4.Range("F20").Value = 215
$ws4.Range("F21").Value = 1567
$ws4.Range("F23").Value = 246
$ws4.Range("F24").Value = 369
$ws4.Range("F25").Value = 6
$ws4.Range("F26").Value = 2
$ws4.Range("F28").Value = 2714
$ws4.Range("F30").Value = 1497
$ws4.Range("F31").Value = 72
$ws4.Range("F32").Value = 69
$ws4.Range("F34").Value = 499
$ws4.Range("F35").Value = 826
$ws4.Range("F36").Value = 1420
$ws4.Range("F40").Value = 1489
$ws4.Range("F41").Value = 804
$ws4.Range("F42").Value = 690
$ws4.Range("F43").Value = 711
$ws4.Range("F44").Value = 910
$ws4.Range("F46").Value = 25
$ws4.Range("F48").Value = 271
$ws4.Range("F49").Value = 40
